$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2025")
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 726.0855929159379
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 20248.23775137067
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 5667.147998863284
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 11319.74679992575
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 40690.0475060112
$ws.Cells.Item(2, 13).Value = 7416.763997874001
$ws.Cells.Item(2, 14).Value = 4971.633266734775
$ws.Cells.Item(2, 15).Value = 4901.543969525117

$ws = $wb.Worksheets.Item("2030")
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 2910.312293597376
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 32194.33633015031
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 5667.147998863284
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 25955.38973957134
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 46019.88203863942
$ws.Cells.Item(2, 13).Value = 12214.33499778223
$ws.Cells.Item(2, 14).Value = 6346.065594672798
$ws.Cells.Item(2, 15).Value = 6820.917285593926

$ws = $wb.Worksheets.Item("2035")
$ws.Cells.Item(2, 1).Value = 2754.31755456332
$ws.Cells.Item(2, 2).Value = 5121.633936870874
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 43659.88035721661
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 5667.147998863284
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 41341.99835421226
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 46019.88203863942
$ws.Cells.Item(2, 13).Value = 16678.1579399765
$ws.Cells.Item(2, 14).Value = 10355.64233412764
$ws.Cells.Item(2, 15).Value = 9956.830726182561

$ws = $wb.Worksheets.Item("2040")
$ws.Cells.Item(2, 1).Value = 2754.31755456332
$ws.Cells.Item(2, 2).Value = 5121.633936870874
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 43659.88035721661
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 5667.147998863284
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 41341.99835421226
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 46019.88203863942
$ws.Cells.Item(2, 13).Value = 16678.1579399765
$ws.Cells.Item(2, 14).Value = 10473.20170275345
$ws.Cells.Item(2, 15).Value = 9956.830726182561

$ws = $wb.Worksheets.Item("2045")
$ws.Cells.Item(2, 1).Value = 5713.151062849596
$ws.Cells.Item(2, 2).Value = 5121.633936870874
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 43659.88035721661
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 5667.147998863284
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 41341.99835421226
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 46019.88203863942
$ws.Cells.Item(2, 13).Value = 16678.1579399765
$ws.Cells.Item(2, 14).Value = 10922.41905690639
$ws.Cells.Item(2, 15).Value = 12033.78963062142

$ws = $wb.Worksheets.Item("2050")
$ws.Cells.Item(2, 1).Value = 5713.151062849596
$ws.Cells.Item(2, 2).Value = 5121.633936870874
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 43659.88035721661
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 5667.147998863284
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 41341.99835421226
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 46019.88203863942
$ws.Cells.Item(2, 13).Value = 16678.1579399765
$ws.Cells.Item(2, 14).Value = 10922.41905690639
$ws.Cells.Item(2, 15).Value = 12033.78963062142
